$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row being appended mirrors the previous row's shape/style; grab the
# previous (last data) row's style so the new cells end up with the same
# default (unstyled) formatting instead of Excel's auto date-detection style.
$prevRow = 61
$newRow = 62

$ws.Cells.Item($newRow, 1).Value = "'2025/10/05"
$ws.Cells.Item($newRow, 1).Style = $ws.Cells.Item($prevRow, 1).Style

$ws.Cells.Item($newRow, 2).Value = "日"
$ws.Cells.Item($newRow, 2).Style = $ws.Cells.Item($prevRow, 2).Style

$ws.Cells.Item($newRow, 3).Value = 4
$ws.Cells.Item($newRow, 3).Style = $ws.Cells.Item($prevRow, 3).Style

$ws.Cells.Item($newRow, 4).Value = 5
$ws.Cells.Item($newRow, 4).Style = $ws.Cells.Item($prevRow, 4).Style
